# Weekly fruit/vegetable price update:
# A new weekly record (row) is inserted for "Primera" quality, dated 44589
# (2022-01-28), pushing all subsequent rows (old 46-88) down by one to
# (47-89), matching the original sheet's row-by-row data shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46; Excel shifts rows 46-88 down to 47-89
# and carries the existing D-column (date) number format onto the new row.
$ws.Rows.Item(46).Insert()

$ws.Cells.Item(46, 1).Value = 1
$ws.Cells.Item(46, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value = 44589
$ws.Cells.Item(46, 5).Value = 15
$ws.Cells.Item(46, 6).Value = 100112036
$ws.Cells.Item(46, 7).Value = "Caigua"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 120
$ws.Cells.Item(46, 11).Value = 11000
$ws.Cells.Item(46, 12).Value = 12000
$ws.Cells.Item(46, 13).Value = 11500
$ws.Cells.Item(46, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 575
$ws.Cells.Item(46, 17).Value = 20
$ws.Cells.Item(46, 18).Value = "Hortaliza"
